$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name in A4 from "DANA SKRUGER" to "JONAS SAMPAIO"
$ws.Range("A4").Value = "JONAS SAMPAIO"

# Update the active selection to A4 (matches the sheetView selection change)
$ws.Range("A4").Select()
